$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round Q2 and R2 to integers
$ws.Range("Q2").Value = 567436
$ws.Range("R2").Value = 6820974

# Clear Z2 and AB2 (Starttid / Sluttid) entirely - delete the cell content
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
